# Projekt Onlineshop FIAEB: Gantt
# Fills in the Gantt-chart day cells (columns L..BM) for the first few
# work packages, un-hides the helper "Kapazitaet" column K, drops the
# now-superfluous "Ressourcen" helper column (BN) together with its
# shared-string values, widens the conditional-formatting range to the
# full, now-contiguous L2:BM17 block, and moves the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gantt")
$ws.Activate()

# --- Gantt bars: mark the active days for rows 2-5 and 7 -------------
$ws.Range("L2:P2").Value = 1
$ws.Range("S3").Value = 1
$ws.Range("T4:W4").Value = 1
$ws.Range("T5:V5").Value = 1
$ws.Range("S7").Value = 1

# --- drop the helper "Ressourcen" values in column BN (rows 2-17) ----
$ws.Range("BN2:BN17").ClearContents()

# --- column K ("Kapazitaet") is no longer hidden ----------------------
$ws.Columns.Item(11).Hidden = $false

# --- widen the "highlight when 1" conditional format to the whole,
#     now-contiguous block L2:BM17 --------------------------------------
$fcs = $ws.Cells.FormatConditions
$fcs.Item(2).ModifyAppliesToRange($ws.Range("L2:BM17"))

# --- move the active selection ----------------------------------------
$excel.Goto($ws.Range("O1"), $true)
$ws.Range("L2").Select()
